# Apply "insurance, claim, debt, investment" column-completion edit (#5)
# for the 債權 (claim) and 債務 (debt) sheets: give them a proper header
# row and the extra legislator/source columns (H:N) that the other sheets
# (土地/建物/汽車/存款/股票) already carry.

$wb = $excel.ActiveWorkbook

function Set-LiteralText($cell, $text) {
    # Force a string like "2012-04-18" to be stored as literal text rather
    # than being auto-converted to a date serial number.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# ---- 債權 (claim) sheet -------------------------------------------------
$ws = $wb.Worksheets.Item("債權")

# Header row (row 1) — was holding stray duplicate data, now proper labels
$ws.Cells.Item(1, 2).Value = "species"
$ws.Cells.Item(1, 3).Value = "owner"
$ws.Cells.Item(1, 4).Value = "debtor"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "register_date"
$ws.Cells.Item(1, 7).Value = "register_reason"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Data row (row 2) — add the missing property_category/category/date/
# legislator_name/legislator_id/source_file/index columns
$ws.Cells.Item(2, 8).Value = "claim"
$ws.Cells.Item(2, 9).Value = "normal"
Set-LiteralText $ws.Cells.Item(2, 10) "2012-04-18"
$ws.Cells.Item(2, 11).Value = "蔡正元"
$ws.Cells.Item(2, 12).Value = 966
$ws.Cells.Item(2, 13).Value = "tmp671f1"
$ws.Cells.Item(2, 14).Value = 118

# ---- 債務 (debt) sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("債務")

# Header row (row 1)
$ws.Cells.Item(1, 2).Value = "species"
$ws.Cells.Item(1, 3).Value = "debtor"
$ws.Cells.Item(1, 4).Value = "owner"
$ws.Cells.Item(1, 5).Value = "total"
$ws.Cells.Item(1, 6).Value = "register_date"
$ws.Cells.Item(1, 7).Value = "register_reason"
$ws.Cells.Item(1, 8).Value = "property_category"
$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 10).Value = "date"
$ws.Cells.Item(1, 11).Value = "legislator_name"
$ws.Cells.Item(1, 12).Value = "legislator_id"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

# Data row (row 2)
$ws.Cells.Item(2, 8).Value = "debt"
$ws.Cells.Item(2, 9).Value = "normal"
Set-LiteralText $ws.Cells.Item(2, 10) "2012-04-18"
$ws.Cells.Item(2, 11).Value = "蔡正元"
$ws.Cells.Item(2, 12).Value = 966
$ws.Cells.Item(2, 13).Value = "tmp671f1"
$ws.Cells.Item(2, 14).Value = 123
